$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C7"  = -13.0064
    "A9"  = -21.81720000000002
    "C12" = -11.2243
    "D15" = -8.693199999999999
    "A18" = -22.15700000000002
    "A20" = -20.97009999999998
    "C26" = -12.82880000000001
    "A27" = -21.73589999999999
    "C27" = -12.7668
    "C29" = -11.3289
    "C37" = -14.65009999999999
    "C38" = -13.5181
    "D38" = -8.820999999999991
    "D44" = -7.601899999999998
    "C51" = -11.6603
    "D51" = -7.466599999999999
    "C55" = -13.8024
    "D57" = -8.0692
    "D63" = -7.764400000000002
    "A69" = -21.92280000000001
    "C69" = -11.6772
    "C70" = -11.7314
    "D70" = -7.783699999999997
    "A76" = -19.86619999999999
    "A82" = -22.16520000000001
    "C83" = -13.98769999999999
    "D99" = -7.806499999999996
    "C102" = -13.5198
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
